$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Aktivitas" (B) / "Batch" (C) columns.
# This shifts old B -> D and old C -> E.
$ws.Range("B1:C3").EntireColumn.Insert()

# Set header row
$ws.Range("B1").Value = "Email"
$ws.Range("C1").Value = "NoTelp"

# Row 2 (Budi Santoso)
$ws.Range("B2").Value = "budi.santoso@email.com"
$ws.Range("C2").Value = "081234567890"

# Row 3 (Siti Aminah)
$ws.Range("B3").Value = "siti.aminah@email.com"
$ws.Range("C3").Value = "081234567891"

# Row 4 (new record - Ahmad Rizki)
$ws.Range("A4").Value = "Ahmad Rizki"
$ws.Range("B4").Value = "ahmad.rizki@email.com"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "Workshop Design Thinking"
$ws.Range("E4").Value = "III"

# Column widths
$ws.Columns.Item(3).ColumnWidth = 15.83203125
$ws.Columns.Item(4).ColumnWidth = 30.83203125
$ws.Columns.Item(5).ColumnWidth = 8.83203125
